$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.097.06'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '1.665.52'
$ws.Range("E3").Value = '  -1.24%  '
$ws.Range("E4").Value = '  -0.89%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.02'
$ws.Range("E5").Value = '  -4.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5160'
$ws.Range("E6").Value = '  -4.94%  '
$ws.Range("E7").Value = '  -0.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2629'
$ws.Range("E8").Value = '  -3.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06202'
$ws.Range("E9").Value = '  -3.81%  '
$ws.Range("E10").Value = '  -4.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07496'
$ws.Range("E11").Value = '  -2.47%  '
$ws.Range("D12").Value = '1.667.61'
$ws.Range("E12").Value = '  -1.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.404'
$ws.Range("E13").Value = '  -2.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5565'
$ws.Range("E14").Value = '  -4.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000007873'
$ws.Range("E15").Value = '  -5.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.45'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").Value = '26.082.68'
$ws.Range("E18").Value = '  -0.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.779'
$ws.Range("E19").Value = '  -3.11%  '
$ws.Range("E20").Value = '  -5.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '185.69'
$ws.Range("E21").Value = '  -2.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.124'
$ws.Range("E22").Value = '  -1.78%  '
$ws.Range("E23").Value = '  -0.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '147.17'
$ws.Range("E24").Value = '  -1.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1235'
$ws.Range("E25").Value = '  -5.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.519'
$ws.Range("E26").Value = '  -4.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.73'
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06234'
$ws.Range("E28").Value = '  -1.15%  '
$ws.Range("E29").Value = '  -3.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.270'
$ws.Range("E30").Value = '  -4.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.469'
$ws.Range("E31").Value = '  -2.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.405'
$ws.Range("E32").Value = '  -5.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.612'
$ws.Range("E33").Value = '  -4.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9910'
$ws.Range("E34").Value = '  -4.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.406'
$ws.Range("E35").Value = '  -0.37%  '
$ws.Range("E36").Value = '  -2.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.704'
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.099'
$ws.Range("E38").Value = '  -2.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01604'
$ws.Range("E39").Value = '  -1.29%  '
$ws.Range("D40").Value = '1.073.29'
$ws.Range("E40").Value = '  -3.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8572'
$ws.Range("E41").Value = '  -2.85%  '
$ws.Range("E42").Value = '  -1.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '98.77'
$ws.Range("E43").Value = '  -2.50%  '
$ws.Range("D44").Value = '1.811.35'
$ws.Range("E44").Value = '  -1.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000110'
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.88'
$ws.Range("E46").Value = '  -2.53%  '
$ws.Range("E47").Value = '  -0.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.939'
$ws.Range("E49").Value = '  -3.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4252'
$ws.Range("E50").Value = '  -1.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.887'
$ws.Range("E51").Value = '  -2.48%  '
